# Generate Report for Handoff
#
# The row for "bda16df1-9ff0-40d0-ab23-1c07c34d60c3.md" has been handed off
# for localization: its status moves from "In Translation" to
# "Ready for handoff", its priority changes from "ht" to "mt", and the
# handoff timestamps are refreshed, on the Overview sheet as well as the
# per-locale (zh-cn / de-de) sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 (bda16df1-...md)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 06:13:48"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 (bda16df1-...md)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-20 06:13:45"

# ---------------------------------------------------------------------
# de-de sheet - row 3 (bda16df1-...md)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-20 06:13:48"

# ---------------------------------------------------------------------
# Widen the Status columns to fit the new, longer "Ready for handoff"
# text (mirrors the column auto-sizing that Excel performs after the
# cell content changes).
# ---------------------------------------------------------------------
$newWidth = 16.333333333333332

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
